# Auto-generated edit script: updates market-price-derived profit columns (H-N)
# across 8 item-category worksheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 272.34286
$ws.Range("I52").Value = 143.5
$ws.Range("K52").Value = 430.5
$ws.Range("M52").Value = -270.5
$ws.Range("H113").Value = 29416372
$ws.Range("I113").Value = 2550.75
$ws.Range("K113").Value = 2550.75
$ws.Range("M113").Value = 703.25
$ws.Range("H128").Value = 89000
$ws.Range("J128").Value = 89000
$ws.Range("L128").Value = 89000
$ws.Range("N128").Value = -98960
$ws.Range("H134").Value = 139333.33
$ws.Range("J134").Value = 139333.33
$ws.Range("L134").Value = 139333.33
$ws.Range("N134").Value = -149473.33
$ws.Range("H138").Value = 6417.1353
$ws.Range("J138").Value = 6516.3335
$ws.Range("L138").Value = 19549.0005
$ws.Range("N138").Value = -29829.0005

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 50043.668
$ws.Range("I44").Value = 16045
$ws.Range("K44").Value = 16045
$ws.Range("M44").Value = -15557
$ws.Range("H61").Value = 33342104
$ws.Range("J61").Value = 45464236
$ws.Range("L61").Value = 45464236
$ws.Range("N61").Value = -45464660
$ws.Range("H74").Value = 49287.453
$ws.Range("I74").Value = 93846.73
$ws.Range("K74").Value = 93846.73
$ws.Range("M74").Value = -92972.73
$ws.Range("H77").Value = 49287.453
$ws.Range("I77").Value = 93846.73
$ws.Range("K77").Value = 469233.65
$ws.Range("M77").Value = -464865.65
$ws.Range("H122").Value = 6679.6
$ws.Range("I122").Value = 6349.5
$ws.Range("K122").Value = 19048.5
$ws.Range("M122").Value = -16598.5
$ws.Range("H132").Value = 6196.543
$ws.Range("I132").Value = 2825.9412
$ws.Range("K132").Value = 8477.8236
$ws.Range("M132").Value = -5947.8236
$ws.Range("H136").Value = 33342104
$ws.Range("J136").Value = 45464236
$ws.Range("L136").Value = 136392708
$ws.Range("N136").Value = -136397808

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5556005
$ws.Range("I22").Value = 6173317
$ws.Range("K22").Value = 6173317
$ws.Range("M22").Value = -6173144
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 4635305
$ws.Range("I134").Value = 8623434
$ws.Range("K134").Value = 25870302
$ws.Range("M134").Value = -25867767

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16686.95
$ws.Range("I31").Value = 5861.6665
$ws.Range("J31").Value = 25544
$ws.Range("K31").Value = 5861.6665
$ws.Range("L31").Value = 25544
$ws.Range("M31").Value = -5566.6665
$ws.Range("N31").Value = -26134
$ws.Range("H34").Value = 16686.95
$ws.Range("I34").Value = 5861.6665
$ws.Range("J34").Value = 25544
$ws.Range("K34").Value = 5861.6665
$ws.Range("L34").Value = 25544
$ws.Range("M34").Value = -5659.6665
$ws.Range("N34").Value = -25948
$ws.Range("H58").Value = 8287.462
$ws.Range("I58").Value = 5112.3076
$ws.Range("J58").Value = 11462.615
$ws.Range("K58").Value = 5112.3076
$ws.Range("L58").Value = 11462.615
$ws.Range("M58").Value = -4909.3076
$ws.Range("N58").Value = -11868.615
$ws.Range("H99").Value = 5425.3687
$ws.Range("I99").Value = 3817.625
$ws.Range("J99").Value = 6594.636
$ws.Range("K99").Value = 3817.625
$ws.Range("L99").Value = 6594.636
$ws.Range("M99").Value = -2319.625
$ws.Range("N99").Value = -9590.636
$ws.Range("H126").Value = 5425.3687
$ws.Range("I126").Value = 3817.625
$ws.Range("J126").Value = 6594.636
$ws.Range("K126").Value = 11452.875
$ws.Range("L126").Value = 19783.908
$ws.Range("M126").Value = -8982.875
$ws.Range("N126").Value = -24723.908
$ws.Range("H134").Value = 7946
$ws.Range("I134").Value = 2907.8333
$ws.Range("J134").Value = 9625.388999999999
$ws.Range("K134").Value = 8723.499899999999
$ws.Range("L134").Value = 28876.167
$ws.Range("M134").Value = -6188.499899999999
$ws.Range("N134").Value = -33946.167
$ws.Range("H136").Value = 8287.462
$ws.Range("I136").Value = 5112.3076
$ws.Range("J136").Value = 11462.615
$ws.Range("K136").Value = 15336.9228
$ws.Range("L136").Value = 34387.845
$ws.Range("M136").Value = -12786.9228
$ws.Range("N136").Value = -39487.845

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 4980.857
$ws.Range("J48").Value = 4980.857
$ws.Range("L48").Value = 14942.571
$ws.Range("N48").Value = -15442.571
$ws.Range("H86").Value = 421.66666
$ws.Range("I86").Value = 386
$ws.Range("K86").Value = 1158
$ws.Range("M86").Value = 28
$ws.Range("H89").Value = 421.66666
$ws.Range("I89").Value = 386
$ws.Range("K89").Value = 3474
$ws.Range("M89").Value = 2454
$ws.Range("H107").Value = 22000452
$ws.Range("J107").Value = 35000390
$ws.Range("L107").Value = 105001170
$ws.Range("N107").Value = -105005010
$ws.Range("H132").Value = 3566.5
$ws.Range("J132").Value = 6299.2
$ws.Range("L132").Value = 56692.8
$ws.Range("N132").Value = -61752.8

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5032.8
$ws.Range("J80").Value = 5555
$ws.Range("L80").Value = 5555
$ws.Range("N80").Value = -7551
$ws.Range("H83").Value = 5032.8
$ws.Range("J83").Value = 5555
$ws.Range("L83").Value = 27775
$ws.Range("N83").Value = -37759
$ws.Range("H93").Value = 39991
$ws.Range("J93").Value = 39991
$ws.Range("L93").Value = 39991
$ws.Range("N93").Value = -43735
$ws.Range("H132").Value = 4848.25
$ws.Range("J132").Value = 8638.083000000001
$ws.Range("L132").Value = 25914.249
$ws.Range("N132").Value = -30974.249
$ws.Range("H134").Value = 96666.664
$ws.Range("J134").Value = 96666.664
$ws.Range("L134").Value = 289999.992
$ws.Range("N134").Value = -295069.992

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5160.4
$ws.Range("I68").Value = 3700.5
$ws.Range("K68").Value = 3700.5
$ws.Range("M68").Value = -2951.5
$ws.Range("H71").Value = 5160.4
$ws.Range("I71").Value = 3700.5
$ws.Range("K71").Value = 18502.5
$ws.Range("M71").Value = -14758.5
$ws.Range("H132").Value = 20009732
$ws.Range("I132").Value = 35719450
$ws.Range("K132").Value = 107158350
$ws.Range("M132").Value = -107155820
$ws.Range("H136").Value = 15289.081
$ws.Range("I136").Value = 8036.5454
$ws.Range("K136").Value = 24109.6362
$ws.Range("M136").Value = -21559.6362

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 201333.33
$ws.Range("I62").Value = 201333.33
$ws.Range("K62").Value = 201333.33
$ws.Range("M62").Value = -200709.33
$ws.Range("H65").Value = 201333.33
$ws.Range("I65").Value = 201333.33
$ws.Range("K65").Value = 1006666.65
$ws.Range("M65").Value = -1003546.65
$ws.Range("H126").Value = 1202.7693
$ws.Range("I126").Value = 1273.6
$ws.Range("J126").Value = 966.6667
$ws.Range("K126").Value = 3820.8
$ws.Range("L126").Value = 2900.0001
$ws.Range("M126").Value = -1350.8
$ws.Range("N126").Value = -7840.0001
$ws.Range("H132").Value = 12558.363
$ws.Range("I132").Value = 9997.154
$ws.Range("J132").Value = 22071.428
$ws.Range("K132").Value = 29991.462
$ws.Range("L132").Value = 66214.284
$ws.Range("M132").Value = -27461.462
$ws.Range("N132").Value = -71274.284
$ws.Range("H136").Value = 22448876
$ws.Range("I136").Value = 41668684
$ws.Range("J136").Value = 483381.44
$ws.Range("K136").Value = 125006052
$ws.Range("L136").Value = 1450144.32
$ws.Range("M136").Value = -125003502
$ws.Range("N136").Value = -1455244.32
